$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B4").Value = 8.775499999999999
$ws.Range("C4").Value = -12.23969999999999
$ws.Range("B6").Value = 6.227799999999998
$ws.Range("D6").Value = -8.322900000000002
$ws.Range("B7").Value = 5.459600000000002
$ws.Range("D7").Value = -8.164899999999999
$ws.Range("B8").Value = 6.166799999999999
$ws.Range("C8").Value = -13.12929999999999
$ws.Range("D8").Value = -8.0509
$ws.Range("C9").Value = -10.2494
$ws.Range("D10").Value = -7.535899999999999
$ws.Range("C12").Value = -10.5743
$ws.Range("D13").Value = -8.335899999999995
$ws.Range("D14").Value = -7.1824
$ws.Range("B16").Value = 4.611799999999998
$ws.Range("D16").Value = -8.159499999999998
$ws.Range("C17").Value = -13.94679999999999
$ws.Range("C18").Value = -11.82489999999999
$ws.Range("C19").Value = -11.12510000000001
$ws.Range("B20").Value = 9.517699999999989
$ws.Range("C20").Value = -12.8341
$ws.Range("B21").Value = 9.193099999999994
$ws.Range("C26").Value = -13.48450000000001
$ws.Range("B28").Value = 6.130600000000004
$ws.Range("B29").Value = 5.323300000000005
$ws.Range("B30").Value = 5.565499999999997
$ws.Range("D30").Value = -6.624999999999999
$ws.Range("C31").Value = -13.3885
$ws.Range("B32").Value = 6.991499999999997
$ws.Range("D37").Value = -7.782499999999996
$ws.Range("C39").Value = -11.355
$ws.Range("B40").Value = 9.245099999999995
$ws.Range("C40").Value = -12.42270000000001
$ws.Range("D40").Value = -8.869099999999994
$ws.Range("C41").Value = -12.47700000000001
$ws.Range("C42").Value = -11.644
$ws.Range("C43").Value = -12.3323
$ws.Range("D44").Value = -6.896600000000007
$ws.Range("B46").Value = 5.547900000000003
$ws.Range("C47").Value = -12.72509999999999
$ws.Range("C48").Value = -12.2211
$ws.Range("B51").Value = 5.620699999999998
$ws.Range("B52").Value = 5.628699999999998
$ws.Range("C54").Value = -12.37710000000001
$ws.Range("B57").Value = 5.180099999999999
$ws.Range("B59").Value = 4.726300000000004
$ws.Range("B62").Value = 5.312799999999998
$ws.Range("C62").Value = -14.20839999999999
$ws.Range("C63").Value = -10.2345
$ws.Range("C64").Value = -10.4746
$ws.Range("B66").Value = 6.3873
$ws.Range("D70").Value = -6.700299999999998
$ws.Range("B73").Value = 8.970499999999998
$ws.Range("B74").Value = 9.003499999999992
$ws.Range("C76").Value = -11.995
$ws.Range("B77").Value = 9.096700000000004
$ws.Range("C81").Value = -14.14719999999999
$ws.Range("C84").Value = -14.11579999999999
$ws.Range("C89").Value = -12.8969
$ws.Range("D89").Value = -8.243999999999993
$ws.Range("D91").Value = -7.669800000000001
$ws.Range("B92").Value = 4.901599999999997
$ws.Range("D93").Value = -6.448499999999993
$ws.Range("C94").Value = -10.66049999999999
$ws.Range("D98").Value = -7.732300000000009
$ws.Range("B100").Value = 5.818199999999998
